$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Drop the trailing period off "NoSQL, bson, MongoDB." and make
#    room for three new space-separated hyperlinks by dropping in
#    unique placeholder tokens we can re-find one at a time. Doing
#    this as a single Find/Replace keeps the edit anchored to the
#    exact existing sentence instead of guessing character offsets.
# -----------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute(
    "NoSQL, bson, MongoDB.",  # FindText
    $false,                   # MatchCase
    $false,                   # MatchWholeWord
    $false,                   # MatchWildcards
    $false,                   # MatchSoundsLike
    $false,                   # MatchAllWordForms
    $true,                    # Forward
    1,                        # Wrap (wdFindContinue)
    $false,                   # Format
    "NoSQL, bson, MongoDB [[MONGO_LINK_1]] [[MONGO_LINK_2]] [[MONGO_LINK_3]]", # ReplaceWith
    2                         # Replace (wdReplaceAll)
) | Out-Null

# -----------------------------------------------------------------
# 2. Turn each placeholder token into its own hyperlink run. Using a
#    fresh Content range + Find for each one (rather than reusing /
#    re-collapsing a previous Range object) keeps the target Range
#    accurate after the prior mutation.
# -----------------------------------------------------------------
$link1 = $d.Content
$link1.Find.Execute("[[MONGO_LINK_1]]") | Out-Null
$d.Hyperlinks.Add(
    $link1,
    "http://github.com/KenAdeniji/WordEngineering/tree/main/IIS/WordEngineering/MongoDB",
    "",
    "",
    "http://github.com/KenAdeniji/WordEngineering/tree/main/IIS/WordEngineering/MongoDB",
    ""
) | Out-Null

$link2 = $d.Content
$link2.Find.Execute("[[MONGO_LINK_2]]") | Out-Null
$d.Hyperlinks.Add(
    $link2,
    "http://github.com/KenAdeniji/WordEngineering/blob/main/InformationInTransit/MongoDB/MongoDBDriverQuickTour.cs",
    "",
    "",
    "http://github.com/KenAdeniji/WordEngineering/blob/main/InformationInTransit/MongoDB/MongoDBDriverQuickTour.cs",
    ""
) | Out-Null

$link3 = $d.Content
$link3.Find.Execute("[[MONGO_LINK_3]]") | Out-Null
$d.Hyperlinks.Add(
    $link3,
    "http://github.com/KenAdeniji/WordEngineering/blob/main/InformationInTransit/MongoDB/MongoDBHelper.cs",
    "",
    "",
    "http://github.com/KenAdeniji/WordEngineering/blob/main/InformationInTransit/MongoDB/MongoDBHelper.cs",
    ""
) | Out-Null

Write-Host "MongoDB hyperlinks inserted"
